# Applies odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "V2"  = 1.58
    "G4"  = 1.23
    "H4"  = 5.6
    "I4"  = 12.5
    "J4"  = 1.62
    "K4"  = 2.6
    "L4"  = 9
    "P4"  = 4.6
    "W4"  = 6.9
    "X4"  = 6.5
    "AD4" = 12
    "AH4" = 28
    "AI4" = 120
    "AJ4" = 40
    "AK4" = 500
    "AU4" = 9
    "AW4" = 11.5
    "AX4" = 65
    "AY4" = 55
    "U5"  = 1.8
    "V5"  = 1.8
    "Q6"  = 2.5
    "R6"  = 1.5
    "U6"  = 2.2
    "V6"  = 1.62
    "X6"  = 7.5
    "Y6"  = 9.5
    "AC6" = 6.5
    "AE6" = 21
    "N10" = 9
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
